$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17206335067749
$ws.Range("B1").Value = 1.354991436004639
$ws.Range("C1").Value = 1.729769706726074
$ws.Range("D1").Value = 3.714582920074463
$ws.Range("E1").Value = 3.691559791564941
